# Insert a new weekly record at row 291 of the "Zapallo italiano" data sheet.
# This pushes every existing record (rows 291-344) down by one row, so the
# former row 344 becomes row 345, and the sheet's used range grows from
# A1:R344 to A1:R345.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 291..344 down by one (Insert defaults to shifting cells down).
$ws.Rows.Item(291).Insert()

# Populate the newly-opened row 291 with the new record's data.
$ws.Cells.Item(291, 1).Value  = 5
$ws.Cells.Item(291, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(291, 3).Value  = "Maule"
$ws.Cells.Item(291, 4).Value  = 44694
$ws.Cells.Item(291, 5).Value  = 7
$ws.Cells.Item(291, 6).Value  = 100112032
$ws.Cells.Item(291, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(291, 8).Value  = "Sin especificar"
$ws.Cells.Item(291, 9).Value  = "Primera"
$ws.Cells.Item(291, 10).Value = 300
$ws.Cells.Item(291, 11).Value = 18000
$ws.Cells.Item(291, 12).Value = 18000
$ws.Cells.Item(291, 13).Value = 18000
$ws.Cells.Item(291, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(291, 15).Value = "Provincia del Elqu" + [char]0x00ED
$ws.Cells.Item(291, 16).Value = 360
$ws.Cells.Item(291, 17).Value = 50
$ws.Cells.Item(291, 18).Value = "Hortaliza"
